$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.224.56"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "1.811.07"
$ws.Range("E3").Value = "  -3.56%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'310.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4207"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.36%  "
$ws.Range("D8").Value = "'0.3560"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.30%  "
$ws.Range("D9").Value = "'0.07122"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.98%  "
$ws.Range("D10").Value = "'0.8506"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("D11").Value = "'20.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.50%  "
$ws.Range("D12").Value = "1.850.53"
$ws.Range("E12").Value = "  -4.26%  "
$ws.Range("D13").Value = "'5.312"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").Value = "'6.390"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").Value = "'0.06863"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "'0.000008774"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.19%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'15.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "27.194.99"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "'5.124"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "'10.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "2.027.63"
$ws.Range("E24").Value = "  -6.97%  "
$ws.Range("D25").Value = "'1.969"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'153.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'18.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("E28").Value = "  -6.85%  "
$ws.Range("D29").Value = "'113.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.50%  "
$ws.Range("E30").Value = "  -9.03%  "
$ws.Range("D31").Value = "'0.08905"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("E32").Value = "  -6.75%  "
$ws.Range("D33").Value = "'2.917"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").Value = "'4.452"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.77%  "
$ws.Range("D35").Value = "'1.110"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.37%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "'1.066"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.66%  "
$ws.Range("D38").Value = "'0.05193"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("D40").Value = "'0.1638"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.4966"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.697"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.75%  "
$ws.Range("D43").Value = "'6.279"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.92%  "
$ws.Range("D44").Value = "'8.179"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.21%  "
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'10.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "'0.06378"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("D49").Value = "'0.4572"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.34%  "
$ws.Range("D50").Value = "'1.594"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("E51").Value = "  -4.11%  "
